$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the F-column (time_taken) timestamps on the "data" sheet ---
$newTimestamps = @(
    "2021-10-05 14:21:32.351484",
    "2021-10-05 14:21:32.351492",
    "2021-10-05 14:21:32.351496",
    "2021-10-05 14:21:32.351498",
    "2021-10-05 14:21:32.351501",
    "2021-10-05 14:21:32.351504",
    "2021-10-05 14:21:32.351506",
    "2021-10-05 14:21:32.351509",
    "2021-10-05 14:21:32.351512",
    "2021-10-05 14:21:32.351514",
    "2021-10-05 14:21:32.351517",
    "2021-10-05 14:21:32.351519",
    "2021-10-05 14:21:32.351522",
    "2021-10-05 14:21:32.351524",
    "2021-10-05 14:21:32.351527",
    "2021-10-05 14:21:32.351529",
    "2021-10-05 14:21:32.351532",
    "2021-10-05 14:21:32.351535",
    "2021-10-05 14:21:32.351537",
    "2021-10-05 14:21:32.351540",
    "2021-10-05 14:21:32.351543",
    "2021-10-05 14:21:32.351546",
    "2021-10-05 14:21:32.351548",
    "2021-10-05 14:21:32.351551",
    "2021-10-05 14:21:32.351553",
    "2021-10-05 14:21:32.351556",
    "2021-10-05 14:21:32.351559",
    "2021-10-05 14:21:32.351562",
    "2021-10-05 14:21:32.351564",
    "2021-10-05 14:21:32.351567",
    "2021-10-05 14:21:32.351569",
    "2021-10-05 14:21:32.351572",
    "2021-10-05 14:21:32.351575",
    "2021-10-05 14:21:32.351578",
    "2021-10-05 14:21:32.351580",
    "2021-10-05 14:21:32.351583",
    "2021-10-05 14:21:32.351586",
    "2021-10-05 14:21:32.351588",
    "2021-10-05 14:21:32.351591",
    "2021-10-05 14:21:32.351593",
    "2021-10-05 14:21:32.351596",
    "2021-10-05 14:21:32.351599",
    "2021-10-05 14:21:32.351601",
    "2021-10-05 14:21:32.351604",
    "2021-10-05 14:21:32.351607",
    "2021-10-05 14:21:32.351609",
    "2021-10-05 14:21:32.351612",
    "2021-10-05 14:21:32.351615",
    "2021-10-05 14:21:32.351617",
    "2021-10-05 14:21:32.351620",
    "2021-10-05 14:21:32.351623"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add the new "metadata" sheet after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (B1:G1) - reuse the same header style as the "data" sheet by
# copying formats only (keeps the same shared style index, no new style
# entries are created in styles.xml).
$metaHeaders = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
$dataSheet.Cells.Item(1, 2).Copy()
for ($c = 0; $c -lt $metaHeaders.Length; $c++) {
    $cell = $metaSheet.Cells.Item(1, $c + 2)
    $cell.PasteSpecial(-4122)
    $cell.Value = $metaHeaders[$c]
}

# A2 uses the same style as the "data" sheet's A-column (index) cells.
$dataSheet.Cells.Item(2, 1).Copy()
$a2 = $metaSheet.Cells.Item(2, 1)
$a2.PasteSpecial(-4122)
$a2.Value = 0

$metaSheet.Cells.Item(2, 2).Value = "Mitochondrial disorder with complex I deficiency"
$metaSheet.Cells.Item(2, 3).Value = 534

# data_version ("1.14") must be stored as text, not a number. Force text
# via a temporary Text number format, then reset the style back to Normal
# (which clears the number format again) so no visible/attached style is
# left on the cell - only the underlying stored value keeps its text type.
$d2 = $metaSheet.Cells.Item(2, 4)
$d2.NumberFormat = "@"
$d2.Value = "1.14"
$d2.Style = "Normal"

$metaSheet.Cells.Item(2, 5).Value = "2021-07-08T10:28:13.875386Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:21:32.348161"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/534/?format=json"

# Keep "data" as the active sheet/tab (matches original activeTab=0)
$dataSheet.Activate()
$excel.CutCopyMode = $false
